$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = "along"

$ws.Range("B2").Value = -0.0359802796359599
$ws.Range("C2").Value = -0.0566583497522772
$ws.Range("D2").Value = -0.0153022095196427
$ws.Range("E2").Value = "cut_aid_in_programTRUE"

$ws.Range("B3").Value = -0.0464175388452988
$ws.Range("C3").Value = -0.0742121605265226
$ws.Range("D3").Value = -0.018622917164075
$ws.Range("E3").Value = "cut_aid_in_programTRUE"

$ws.Range("B4").Value = 0.0153185378491868
$ws.Range("C4").Value = -0.053732438640918
$ws.Range("D4").Value = 0.0843695143392915
$ws.Range("E4").Value = "cut_aid_in_programTRUE"

$ws.Range("B5").Value = -0.0694482909577571
$ws.Range("C5").Value = -0.127593399095205
$ws.Range("D5").Value = -0.0113031828203091
$ws.Range("E5").Value = "cut_aid_in_programTRUE"

$ws.Range("B6").Value = -0.126411963213932
$ws.Range("C6").Value = -0.197866304239028
$ws.Range("D6").Value = -0.0549576221888362
$ws.Range("E6").Value = "cut_aid_in_programTRUE"

$ws.Range("B7").Value = -0.0542406068491518
$ws.Range("C7").Value = -0.13574855367973
$ws.Range("D7").Value = 0.0272673399814261
$ws.Range("E7").Value = "cut_aid_in_programTRUE"

$ws.Range("B8").Value = -0.10179630030163
$ws.Range("C8").Value = -0.17705687228679
$ws.Range("D8").Value = -0.0265357283164696
$ws.Range("E8").Value = "cut_aid_in_programTRUE"

$ws.Range("B9").Value = 0.0108339402477596
$ws.Range("C9").Value = -0.0556931917155847
$ws.Range("D9").Value = 0.0773610722111039
$ws.Range("E9").Value = "cut_aid_in_programTRUE"

$ws.Range("B10").Value = 0.0837839851033947
$ws.Range("C10").Value = -0.0216366707883592
$ws.Range("D10").Value = 0.189204640995149
$ws.Range("E10").Value = "cut_aid_in_programTRUE"

$ws.Range("B11").Value = -0.0148317041662747
$ws.Range("C11").Value = -0.0493904940584233
$ws.Range("D11").Value = 0.0197270857258738
$ws.Range("E11").Value = "cut_aid_in_programTRUE"

$ws.Range("B12").Value = -0.0384451845061834
$ws.Range("C12").Value = -0.0987851727871903
$ws.Range("D12").Value = 0.0218948037748236
$ws.Range("E12").Value = "cut_aid_in_programTRUE"
